# Add a new "2022-Q3" quarter sheet to the workbook, positioned right after
# the "总计" (summary) sheet and before the existing "2022-Q2" sheet, then
# update the "总计" sheet with a new row summarizing the 2022-Q3 data.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Update "总计" (summary) sheet: insert 2022-Q3 as the new first data
#    row and shift the existing quarters down by one row.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# Copy row 5's formatting down into the freshly used row 6 so the new
# bottom row keeps the same look (bold centered style used by column A).
$summary.Range("A5").Copy()
$summary.Range("A6").PasteSpecial(-4122)

$summary.Range("A6").Value = 4
$summary.Range("B6").Value = "2021-Q2"
$summary.Range("C6").Value = 5
$summary.Range("D6").Value = 3.83

$summary.Range("A5").Value = 3
$summary.Range("B5").Value = "2021-Q3"
$summary.Range("C5").Value = 1
$summary.Range("D5").Value = 1.92

$summary.Range("A4").Value = 2
$summary.Range("B4").Value = "2022-Q1"
$summary.Range("C4").Value = 4
$summary.Range("D4").Value = 2.31

$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2022-Q2"
$summary.Range("C3").Value = 4
$summary.Range("D3").Value = 2.71

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 4
$summary.Range("D2").Value = 2.54

# ---------------------------------------------------------------------
# 2) Insert the new "2022-Q3" worksheet right before "2022-Q2" by
#    duplicating the "2022-Q2" sheet (same columns/layout/styles) and
#    then updating the figures that differ for the new quarter. This
#    keeps formatting (borders, bold header, text-typed numeric-looking
#    columns, sheet page setup, …) identical to the existing quarter
#    sheets, matching how this recurring report is actually produced.
# ---------------------------------------------------------------------
$q2Sheet = $wb.Worksheets.Item("2022-Q2")
$q2Sheet.Copy($q2Sheet)
$q3Sheet = $wb.Worksheets.Item(2)
$q3Sheet.Name = "2022-Q3"

# Figures in columns D-G are stored as text (matching the source
# workbook's numeric-looking-but-textual cells), so force text number
# formatting before overwriting them - otherwise Excel auto-converts
# the numeric-looking strings to real numbers.
$q3Sheet.Range("D2:G5").NumberFormat = "@"

# Row 2 - 968029 (name/code unchanged, only the figures moved)
$q3Sheet.Range("D2").Value = "27.03"
$q3Sheet.Range("E2").Value = "99.07"
$q3Sheet.Range("F2").Value = "7.90"
$q3Sheet.Range("G2").Value = "2.1354"
$q3Sheet.Range("H2").Value = 2

# Row 3 - 009562
$q3Sheet.Range("C3").Value = "工银全球股票（QDII）美元"
$q3Sheet.Range("D3").Value = "5.89"
$q3Sheet.Range("E3").Value = "93.72"
$q3Sheet.Range("F3").Value = "2.30"
$q3Sheet.Range("G3").Value = "0.1355"
$q3Sheet.Range("H3").Value = 5

# Row 4 - 009563
$q3Sheet.Range("C4").Value = "工银全球股票（QDII）港币"
$q3Sheet.Range("D4").Value = "5.89"
$q3Sheet.Range("E4").Value = "93.72"
$q3Sheet.Range("F4").Value = "2.30"
$q3Sheet.Range("G4").Value = "0.1355"
$q3Sheet.Range("H4").Value = 5

# Row 5 - 486001
$q3Sheet.Range("C5").Value = "工银瑞信中国机会全球配置股票（QDII）人民币"
$q3Sheet.Range("D5").Value = "5.89"
$q3Sheet.Range("E5").Value = "93.72"
$q3Sheet.Range("F5").Value = "2.30"
$q3Sheet.Range("G5").Value = "0.1355"
$q3Sheet.Range("H5").Value = 5

# Restore the originally active sheet/selection (the last quarter tab was
# the selected one before this edit, and stays so afterwards).
$lastSheet = $wb.Worksheets.Item("2021-Q2")
$lastSheet.Activate()
$lastSheet.Range("A1").Select()
